$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the three target paragraphs (Graph1 - "Memory vs Problem Size"
# block): "Basic:", "Efficient:" and "Explanation: ".
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$basicIdx = 0
$efficientIdx = 0
$explanationIdx = 0
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd([char]13)
    if ($basicIdx -eq 0 -and $t -eq "Basic:") {
        $basicIdx = $i
    }
    if ($basicIdx -ne 0 -and $efficientIdx -eq 0 -and $t -eq "Efficient:") {
        $efficientIdx = $i
    }
    if ($efficientIdx -ne 0 -and $explanationIdx -eq 0 -and $t -eq "Explanation: ") {
        $explanationIdx = $i
        break
    }
}

$nsAttr = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) "Basic:" paragraph gains a second run: " Exponential"
# ---------------------------------------------------------------------
$pBasic = $d.Paragraphs.Item($basicIdx)
$rBasic = $pBasic.Range
$xmlBasic = '<w:p' + $nsAttr + '><w:r><w:t>Basic:</w:t></w:r><w:r><w:t xml:space="preserve"> Exponential</w:t></w:r></w:p>'
$rBasic.InsertXML($xmlBasic)
Write-Host "updated Basic paragraph"

# ---------------------------------------------------------------------
# 2) "Efficient:" paragraph gains a second run: " Linear"
# ---------------------------------------------------------------------
$pEfficient = $d.Paragraphs.Item($efficientIdx)
$rEfficient = $pEfficient.Range
$xmlEfficient = '<w:p' + $nsAttr + '><w:r><w:t>Efficient:</w:t></w:r><w:r><w:t xml:space="preserve"> Linear</w:t></w:r></w:p>'
$rEfficient.InsertXML($xmlEfficient)
Write-Host "updated Efficient paragraph"

# ---------------------------------------------------------------------
# 3) A brand-new justified paragraph is added right after
#    "Explanation: ", explaining the memory growth.
# ---------------------------------------------------------------------
$pExplanation = $d.Paragraphs.Item($explanationIdx)
$rExplanation = $pExplanation.Range

$sentence1 = "As we can see from the graph, the basic dynamic programming algorithm takes exponentially more memory as the problem size is increased because it requires to create a memorization table of size m*n whereas the efficient algorithm which uses a divide and conquer approach with dynamic programming, we will use"
$sentence2 = " only 2 * max("
$sentence3 = "m,n"
$sentence4 = ") space which explains the linear memory growth."

$xmlExplanation = '<w:p' + $nsAttr + '><w:r><w:t xml:space="preserve">Explanation: </w:t></w:r></w:p>' + `
    '<w:p' + $nsAttr + '><w:pPr><w:jc w:val="both"/></w:pPr>' + `
    '<w:r><w:t>' + $sentence1 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $sentence2 + '</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>' + $sentence3 + '</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>' + $sentence4 + '</w:t></w:r>' + `
    '</w:p>'

$rExplanation.InsertXML($xmlExplanation)
Write-Host "inserted memory explanation paragraph"
